# The research-question slide ("Hypotheses") restates both the null and the
# alternative hypothesis. The author tightened the wording by dropping the
# redundant word "significant" before "difference" in both sentences:
#   "There is no significant difference ..." -> "There is no difference ..."
#   "There is a significant difference ..."  -> "There is a difference ..."
#
# Slide 6 ("Hypotheses") holds this text in shape "Title 8" (shape index 3).

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(6)
$shp = $s.Shapes.Item(3)          # "Title 8"
$tr  = $shp.TextFrame.TextRange

$old = "significant difference "
$new = "difference "

# First sentence ("... There is no significant difference ...").
$full = $tr.Text
$idx  = $full.IndexOf($old)
if ($idx -ge 0) {
    $tr.Characters($idx + 1, $old.Length).Text = $new
}

# Second sentence ("... There is a significant difference ..."). Search
# past the first sentence so this lands on the *second* occurrence.
$full       = $tr.Text
$searchFrom = $full.IndexOf("Alternative hypothesis")
$idx        = $full.IndexOf($old, $searchFrom)
if ($idx -ge 0) {
    $tr.Characters($idx + 1, $old.Length).Text = $new
}
